$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2017-02-28 06:55:56"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2017-02-28 06:56:11"
